# Add a new attribute row "hasOnlyNovelOmics" to the "attributes" sheet,
# inserted immediately after the existing "patch" row (row 8), pushing all
# subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attributes")

# Insert a new row before row 9 (current row that holds "organisation"),
# shifting it (and everything after it) down by one.
$ws.Rows.Item(9).Insert()

$ws.Cells.Item(9, 1).Value = "rd3_overview"
$ws.Cells.Item(9, 2).Value = "hasOnlyNovelOmics"
$ws.Cells.Item(9, 3).Value = "If true, this subject does not have data in any data freeze (1,2,3,etc.)"
$ws.Cells.Item(9, 5).Value = "bool"
$ws.Cells.Item(9, 6).Value = $false
$ws.Cells.Item(9, 7).Value = $false
$ws.Cells.Item(9, 8).Value = $true
$ws.Cells.Item(9, 9).Value = $false
